$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 12: average of the k column (J) just under the data table.
# ---------------------------------------------------------------------------
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11

# ---------------------------------------------------------------------------
# Rows 14-17: summary statistics block (labels in A, formulas in B).
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108
$ws.Rows.Item(14).RowHeight = 15.6

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"
$ws.Range("B15").Font.Bold = $true
$ws.Range("B15").Font.Size = 12
$ws.Range("B15").VerticalAlignment = -4108
$ws.Rows.Item(15).RowHeight = 15.6

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"
$ws.Range("B16").Font.Bold = $true
$ws.Range("B16").Font.Size = 12
$ws.Range("B16").VerticalAlignment = -4108
$ws.Rows.Item(16).RowHeight = 15.6

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"
$ws.Range("B17").Font.Bold = $true
$ws.Range("B17").Font.Size = 12
$ws.Range("B17").VerticalAlignment = -4108
$ws.Rows.Item(17).RowHeight = 15.6

# ---------------------------------------------------------------------------
# Page setup: paper size A4 / portrait orientation.
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Leave the final selection on J12, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("J12").Select()

Write-Output "done"
